$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date/time column (A) gets a date-time number format (maps to builtin numFmtId 22)
# and a wider column to fit it.
$ws.Range("A1:A5").NumberFormat = "m/d/yy h:mm"
$ws.Columns.Item(1).ColumnWidth = 14.08333333333333

# New data rows (2-5), pulled from the SAVE bag log.
$rows = @(
    @(42606.574745370373, -52, 38, 60, 0, 100, 824, 1785, 215, 16, 25, 0, 4, "Bag"),
    @(42606.575821759259, -26, 53, 46, 0, 100, 2404, 2409, 308, 29, 25, 0, 4, "Bag"),
    @(42606.580879629626, -26, 53, 46, 0, 100, 1076, 2417, 308, 29, 25, 0, 4, "Bag"),
    @(42606.581273148149, 0, 50, 48, 50, 50, 2661, 2888, 368, 35, 34, 1, 1, "Bag")
)

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $r).Value = $row[$c]
    }
}
